# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a failed
# handback transform for the "251f6199-...md" record:
#   - Status changes from "Ready for handoff" to "Handback transform failed"
#     on the Overview sheet (zh-cn / de-de columns) and on each language
#     sheet's "Status" column.
#   - The "Error Detail" column (P) on the zh-cn and de-de sheets is filled
#     in with an explanation of the handback filename mismatch, and that
#     column is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: qdrzjs2p.4y2 is different with handoff file name: 251f6199-1781-42d6-aec8-2e290a3b059e.1163eeae539155818cd79bcf44e570bad0ef45d1.zh-cn."
# raw OOXML column width of 40 corresponds to a COM ColumnWidth of 40 - 5/6
$zhcn.Range("P1").ColumnWidth = 39.166666666666664

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: qdrzjs2p.4y2 is different with handoff file name: 251f6199-1781-42d6-aec8-2e290a3b059e.1163eeae539155818cd79bcf44e570bad0ef45d1.de-de."
$dede.Range("P1").ColumnWidth = 39.166666666666664
